# Apply the updated crypto price/volume snapshot (GitHub Actions data refresh).
# Rows 16/17, 28/29, 40/41, 46/47 were also re-sorted (coin name + link swapped
# with their neighbour), so B/C are rewritten there too.
#
# A handful of D-column price values are plain decimals (one ".") which Excel
# would otherwise coerce to Number on assignment (dropping e.g. trailing zeros:
# "15.00" -> 15, "0.3620" -> 0.362). Those are written with a leading quote so
# they stay text, matching the source inlineStr cells exactly. Multi-dot values
# (e.g. "27.102.40") and the padded percent strings in column E are never parsed
# as numbers by Excel, so they need no such prefix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.102.40'
$ws.Range("E2").Value = '  -1.82%  '

# Row 3
$ws.Range("D3").Value = '1.804.39'
$ws.Range("E3").Value = '  -1.84%  '

# Row 4
$ws.Range("D4").Value = '''1.009'
$ws.Range("E4").Value = '  +0.70%  '

# Row 5
$ws.Range("D5").Value = '''313.66'
$ws.Range("E5").Value = '  -0.47%  '

# Row 6
$ws.Range("D6").Value = '''1.008'
$ws.Range("E6").Value = '  +0.70%  '

# Row 7
$ws.Range("D7").Value = '''0.4242'
$ws.Range("E7").Value = '  -1.27%  '

# Row 8
$ws.Range("D8").Value = '''0.3620'
$ws.Range("E8").Value = '  -1.50%  '

# Row 9
$ws.Range("D9").Value = '''0.07206'
$ws.Range("E9").Value = '  -1.08%  '

# Row 10
$ws.Range("D10").Value = '''0.8451'
$ws.Range("E10").Value = '  -3.06%  '

# Row 11
$ws.Range("D11").Value = '''20.30'
$ws.Range("E11").Value = '  -2.12%  '

# Row 12
$ws.Range("D12").Value = '1.837.55'
$ws.Range("E12").Value = '  -0.34%  '

# Row 13
$ws.Range("D13").Value = '''5.294'
$ws.Range("E13").Value = '  -2.72%  '

# Row 14
$ws.Range("D14").Value = '''6.396'
$ws.Range("E14").Value = '  -2.39%  '

# Row 15
$ws.Range("D15").Value = '''0.06841'
$ws.Range("E15").Value = '  -1.59%  '

# Row 16
$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").Value = '''1.011'
$ws.Range("E16").Value = '  +0.78%  '

# Row 17
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").Value = '''81.11'
$ws.Range("E17").Value = '  +0.71%  '

# Row 18
$ws.Range("D18").Value = '''0.000008748'
$ws.Range("E18").Value = '  -2.34%  '

# Row 19
$ws.Range("D19").Value = '''1.008'
$ws.Range("E19").Value = '  +0.60%  '

# Row 20
$ws.Range("D20").Value = '''15.00'
$ws.Range("E20").Value = '  -3.25%  '

# Row 21
$ws.Range("D21").Value = '27.117.89'
$ws.Range("E21").Value = '  -1.83%  '

# Row 22
$ws.Range("D22").Value = '''5.070'
$ws.Range("E22").Value = '  -1.90%  '

# Row 23
$ws.Range("D23").Value = '''11.12'
$ws.Range("E23").Value = '  +1.97%  '

# Row 24
$ws.Range("D24").Value = '2.022.56'
$ws.Range("E24").Value = '  -2.54%  '

# Row 25
$ws.Range("D25").Value = '''1.963'
$ws.Range("E25").Value = '  -1.03%  '

# Row 26
$ws.Range("D26").Value = '''153.07'
$ws.Range("E26").Value = '  -0.82%  '

# Row 27
$ws.Range("D27").Value = '''18.26'
$ws.Range("E27").Value = '  -3.18%  '

# Row 28
$ws.Range("B28").Value = 'BitcoinCash'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D28").Value = '''115.29'
$ws.Range("E28").Value = '  +0.45%  '

# Row 29
$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").Value = '''5.016'
$ws.Range("E29").Value = '  -4.22%  '

# Row 30
$ws.Range("D30").Value = '''1.633'
$ws.Range("E30").Value = '  -11.60%  '

# Row 31
$ws.Range("D31").Value = '''0.08939'
$ws.Range("E31").Value = '  +0.47%  '

# Row 32
$ws.Range("D32").Value = '''0.7289'
$ws.Range("E32").Value = '  -5.31%  '

# Row 33
$ws.Range("D33").Value = '''2.849'
$ws.Range("E33").Value = '  -3.78%  '

# Row 34
$ws.Range("D34").Value = '''4.344'
$ws.Range("E34").Value = '  -4.86%  '

# Row 35
$ws.Range("D35").Value = '''1.097'
$ws.Range("E35").Value = '  -4.53%  '

# Row 36
$ws.Range("D36").Value = '''1.009'
$ws.Range("E36").Value = '  +0.72%  '

# Row 37
$ws.Range("D37").Value = '''1.090'
$ws.Range("E37").Value = '  -0.81%  '

# Row 38
$ws.Range("D38").Value = '''0.01910'
$ws.Range("E38").Value = '  -2.06%  '

# Row 39
$ws.Range("D39").Value = '''0.05118'
$ws.Range("E39").Value = '  -4.01%  '

# Row 40
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = '''0.4975'
$ws.Range("E40").Value = '  -2.71%  '

# Row 41
$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").Value = '''0.1623'
$ws.Range("E41").Value = '  -3.23%  '

# Row 42
$ws.Range("D42").Value = '''2.619'
$ws.Range("E42").Value = '  -6.97%  '

# Row 43
$ws.Range("D43").Value = '''6.002'
$ws.Range("E43").Value = '  -9.59%  '

# Row 44
$ws.Range("D44").Value = '''8.075'
$ws.Range("E44").Value = '  -5.02%  '

# Row 45
$ws.Range("D45").Value = '''10.25'
$ws.Range("E45").Value = '  -3.01%  '

# Row 46
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").Value = '''1.008'
$ws.Range("E46").Value = '  +0.76%  '

# Row 47
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").Value = '''104.71'
$ws.Range("E47").Value = '  -1.60%  '

# Row 48
$ws.Range("D48").Value = '''0.06321'
$ws.Range("E48").Value = '  -2.74%  '

# Row 49
$ws.Range("D49").Value = '''0.4520'
$ws.Range("E49").Value = '  -4.15%  '

# Row 50
$ws.Range("D50").Value = '''1.591'
$ws.Range("E50").Value = '  -2.30%  '

# Row 51
$ws.Range("D51").Value = '''1.717'
$ws.Range("E51").Value = '  -3.49%  '
